$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Navkaran Singh Sidhu"

# Preconditions column (E) for rows 7-11: "None"
$ws.Range("E7:E11").Value = "None"

# Method Inputs column (F) for rows 9-11: quoted Green inputs
$ws.Range("F9:F11").Value = """Green"",9,10,11"

# Expected Result for row 7
$ws.Range("G7").Value = "Attribute set to the argument values."

# Method Inputs column (F) for rows 7-8: unquoted Green inputs
$ws.Range("F7:F8").Value = "Green,9,10,11"

# Preconditions column (E) for rows 12-14: Triangle(...) call
$ws.Range("E12:E14").Value = "Triangle(""Green"",9,10,11)"

# Expected Result column (G) for rows 8-11: ValueError
$ws.Range("G8:G11").Value = "ValueError"

# Expected Result for row 12
$ws.Range("G12").Value = "The shape color is Green"

# Method Inputs column (F) for rows 12-15: None
$ws.Range("F12:F15").Value = "None"

# Row 15 extra columns
$ws.Range("C15").Value = "Triangle inequality theorm"
$ws.Range("D15").Value = "Exception raised when sum of two sides is smaller then the third side."
$ws.Range("E15").Value = "Triangle(""Green"",1,1,11)"
$ws.Range("G15").Value = "ValueError"

# Numeric expected results
$ws.Range("G13").Value = 42.43
$ws.Range("G14").Value = 30

# Update selection to match the author's final cursor position
$ws.Range("E15").Select()
